# Kilimanjaro_Weekly_Scoreboard.xlsx - "Add files via upload"
# Appends 7 new weekly workout rows (week 7, 2024-07-25) to the bottom of
# the Sheet1 log table, then updates the view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new rows to append, in the same column order as the existing table:
# A Participant | B Date | C Workout Type | D Total Duration |
# E Total Distance | F Total Elevation | G Zone1 | H Zone2 | I Zone3 |
# J Zone4 | K Zone5 | L Workout Level | M Week
$newRows = @(
  @{ Row=280; A="Phil";   B=45498; C="Run";     D=7;  E=0.63; F=73;  G=1;  H=2;  I=0; J=0; K=0; L="Sauntering Hippo"; M=7 },
  @{ Row=281; A="Phil";   B=45498; C="Workout"; D=56; E=0;    F=0;   G=26; H=28; I=2; J=0; K=0; L="Sauntering Hippo"; M=7 },
  @{ Row=282; A="Phil";   B=45498; C="Run";     D=13; E=1.56; F=117; G=0;  H=2;  I=6; J=4; K=0; L="Sauntering Hippo"; M=7 },
  @{ Row=283; A="Eric";   B=45498; C="Workout"; D=92; E=0;    F=0;   G=58; H=31; I=4; J=0; K=0; L="Sauntering Hippo"; M=7 },
  @{ Row=284; A="Phil";   B=45498; C="Walk";    D=18; E=1.04; F=44;  G=18; H=0;  I=0; J=0; K=0; L="Sauntering Hippo"; M=7 },
  @{ Row=285; A="Matt";   B=45498; C="Workout"; D=50; E=0;    F=0;   G=28; H=20; I=3; J=0; K=0; L="Agile Antelope";   M=7 },
  @{ Row=286; A="Steven"; B=45498; C="Walk";    D=38; E=2.04; F=66;  G=38; H=0;  I=0; J=0; K=0; L="Brave Leopard";    M=7 }
)

# The existing last data row (279) carries the date-formatted style (s="1")
# on column B; copy it down onto each new B cell so the new date values
# pick up the same number format instead of minting a new style.
$lastRow = 279

foreach ($r in $newRows) {
  $row = $r.Row

  $ws.Range("B$lastRow").Copy($ws.Range("B$row"))

  $ws.Range("A$row").Value = $r.A
  $ws.Range("B$row").Value = $r.B
  $ws.Range("C$row").Value = $r.C
  $ws.Range("D$row").Value = $r.D
  $ws.Range("E$row").Value = $r.E
  $ws.Range("F$row").Value = $r.F
  $ws.Range("G$row").Value = $r.G
  $ws.Range("H$row").Value = $r.H
  $ws.Range("I$row").Value = $r.I
  $ws.Range("J$row").Value = $r.J
  $ws.Range("K$row").Value = $r.K
  $ws.Range("L$row").Value = $r.L
  $ws.Range("M$row").Value = $r.M
}

# Match the saved view state: scrolled down a few rows further and the
# active cell/selection parked past the new data.
$ws.Activate()
$ws.Range("P276").Select()
